$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. "300.46"); Excel
# auto-converts such text to a Number on assignment unless the cell is
# pre-formatted as Text. Force Text format for the whole price column,
# write the values, then reset the style so no stray number format sticks.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '43.068.02'
$ws.Range("E2").Value = '  +0.03%  '
$ws.Range("D3").Value = '2.303.51'
$ws.Range("E3").Value = '  +0.07%  '
$ws.Range("D5").Value = '300.46'
$ws.Range("E5").Value = '  -0.12%  '
$ws.Range("D6").Value = '98.17'
$ws.Range("E6").Value = '  -1.55%  '
$ws.Range("E7").Value = '  +2.83%  '
$ws.Range("D9").Value = '0.517'
$ws.Range("E9").Value = '  +1.19%  '
$ws.Range("D10").Value = '36.13'
$ws.Range("E10").Value = '  -0.54%  '
$ws.Range("D11").Value = '0.0793'
$ws.Range("E11").Value = '  +0.27%  '
$ws.Range("E12").Value = '  +0.78%  '
$ws.Range("E13").Value = '  -2.54%  '
$ws.Range("D14").Value = '6.89'
$ws.Range("D15").Value = '2.659.18'
$ws.Range("E15").Value = '  -0.05%  '
$ws.Range("D16").Value = '2.346.84'
$ws.Range("E16").Value = '  +2.00%  '
$ws.Range("D17").Value = '0.789'
$ws.Range("E17").Value = '  -1.37%  '
$ws.Range("D18").Value = '42.927.12'
$ws.Range("D19").Value = '12.80'
$ws.Range("E19").Value = '  +0.72%  '
$ws.Range("D20").Value = '0.0₃0912'
$ws.Range("E20").Value = '  +0.78%  '
$ws.Range("D21").Value = '6.15'
$ws.Range("E21").Value = '  +0.23%  '
$ws.Range("D22").Value = '68.32'
$ws.Range("E22").Value = '  +0.59%  '
$ws.Range("D23").Value = '238.07'
$ws.Range("E23").Value = '  +1.03%  '
$ws.Range("E24").Value = '  -1.01%  '
$ws.Range("E25").Value = '  -0.46%  '
$ws.Range("D26").Value = '2.43'
$ws.Range("E26").Value = '  -0.61%  '
$ws.Range("E27").Value = '  -0.14%  '
$ws.Range("D28").Value = '25.04'
$ws.Range("E28").Value = '  +0.28%  '
$ws.Range("B29").Value = 'Monero'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D29").Value = '164.06'
$ws.Range("E29").Value = '  -2.34%  '
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").Value = '2.04'
$ws.Range("E30").Value = '  -13.08%  '
$ws.Range("B31").Value = 'Cosmos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D31").Value = '9.15'
$ws.Range("E31").Value = '  +0.05%  '
$ws.Range("D32").Value = '33.13'
$ws.Range("E32").Value = '  -4.26%  '
$ws.Range("D33").Value = '0.999'
$ws.Range("D34").Value = '5.12'
$ws.Range("E34").Value = '  +1.75%  '
$ws.Range("B35").Value = 'RenderToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D35").Value = '4.79'
$ws.Range("E35").Value = '  +4.03%  '
$ws.Range("B36").Value = 'Celestia'
$ws.Range("C36").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D36").Value = '18.10'
$ws.Range("E36").Value = '  +2.50%  '
$ws.Range("E37").Value = '  +0.35%  '
$ws.Range("D38").Value = '0.0697'
$ws.Range("E38").Value = '  +0.92%  '
$ws.Range("E39").Value = '  +1.05%  '
$ws.Range("E40").Value = '  -0.43%  '
$ws.Range("E41").Value = '  -1.02%  '
$ws.Range("E42").Value = '  +1.20%  '
$ws.Range("D43").Value = '2.021.09'
$ws.Range("E43").Value = '  +2.07%  '
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").Value = '0.0287'
$ws.Range("E44").Value = '  -1.43%  '
$ws.Range("B45").Value = 'ApeXProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D45").Value = '2.24'
$ws.Range("E45").Value = '  -2.32%  '
$ws.Range("D46").Value = '10.39'
$ws.Range("E46").Value = '  +2.18%  '
$ws.Range("D47").Value = '17.53'
$ws.Range("E47").Value = '  +0.02%  '
$ws.Range("D48").Value = '2.84'
$ws.Range("E48").Value = '  -2.41%  '
$ws.Range("D49").Value = '54.37'
$ws.Range("E49").Value = '  -1.97%  '
$ws.Range("D50").Value = '2.525.20'
$ws.Range("E50").Value = '  +0.05%  '
$ws.Range("D51").Value = '1.54'
$ws.Range("E51").Value = '  -1.30%  '

# Restore default ("Normal") style on the price column so no extra
# cell formatting is left behind now that the text has been written.
$priceRange.Style = "Normal"
